$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.660.26"
$ws.Range("E2").Value = "  +3.05%  "
$ws.Range("D3").Value = "2.037.63"
$ws.Range("E3").Value = "  +7.50%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.28"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.661"
$ws.Range("E6").Value = "  -4.56%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.98"
$ws.Range("E8").Value = "  +4.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "60.18"
$ws.Range("E9").Value = "  +6.17%  "
$ws.Range("E10").Value = "  +0.87%  "
$ws.Range("E11").Value = "  -4.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0984"
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.50"
$ws.Range("E13").Value = "  -0.60%  "
$ws.Range("D14").Value = "2.337.97"
$ws.Range("E14").Value = "  +7.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.804"
$ws.Range("E15").Value = "  +1.33%  "
$ws.Range("D16").Value = "2.029.31"
$ws.Range("E16").Value = "  +6.89%  "
$ws.Range("E17").Value = "  -3.07%  "
$ws.Range("D18").Value = "36.649.72"
$ws.Range("E18").Value = "  +2.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.05"
$ws.Range("E19").Value = "  -3.42%  "
$ws.Range("E20").Value = "  -2.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "236.79"
$ws.Range("E21").Value = "  -4.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.58"
$ws.Range("E22").Value = "  -3.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.88"
$ws.Range("E23").Value = "  -5.95%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.46"
$ws.Range("E25").Value = "  -8.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.72"
$ws.Range("E26").Value = "  +1.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.76"
$ws.Range("E27").Value = "  +0.68%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.78"
$ws.Range("E28").Value = "  +7.62%  "
$ws.Range("E29").Value = "  -9.69%  "
$ws.Range("E30").Value = "  -5.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "21.56"
$ws.Range("E31").Value = "  +48.89%  "
$ws.Range("E32").Value = "  -1.68%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0578"
$ws.Range("E33").Value = "  -5.05%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0879"
$ws.Range("E35").Value = "  +17.92%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.87"
$ws.Range("E36").Value = "  +1.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.23"
$ws.Range("E37").Value = "  +14.42%  "
$ws.Range("E39").Value = "  +0.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.33"
$ws.Range("E40").Value = "  -10.60%  "
$ws.Range("E41").Value = "  -6.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "95.77"
$ws.Range("E42").Value = "  -3.59%  "
$ws.Range("E43").Value = "  +2.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.77"
$ws.Range("E44").Value = "  +15.53%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.89"
$ws.Range("E45").Value = "  -6.30%  "
$ws.Range("D46").Value = "1.307.88"
$ws.Range("E46").Value = "  -0.47%  "
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.79"
$ws.Range("E48").Value = "  +1.93%  "
$ws.Range("D49").Value = "2.231.23"
$ws.Range("E49").Value = "  +7.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.20"
$ws.Range("E50").Value = "  -6.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.81"
$ws.Range("E51").Value = "  +14.42%  "
